# Re-pull data, push all data, mean calculation.
# Update the "dSF" (F) column values to reflect the repulled/recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -5
    5  = -4
    6  = -3
    8  = -5
    9  = -6
    10 = -7
    11 = -1
    12 = -2
    13 = -1
    16 = -7
    17 = -6
    18 = -8
    19 = -2
    20 = -2
    21 = 2
    22 = 2
    24 = 4
    25 = -7
    26 = -1
    27 = -3
    28 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
